# Update workbook data/epexspot_prices.xlsx
# 1. "Prix Spot" sheet: insert a new date column ("15-dec") right before the
#    "01-oct." column (currently column EL), shifting all the "01-oct."..
#    "31-oct." columns one place to the right. The new column gets the
#    header "15-dec" and a "-" placeholder in every data row (no price yet).
# 2. "Gaz" sheet: append two new rows for 2025-12-13 and 2025-12-14.
# 3. "CO2" sheet: append two new rows for 2025-12-13 and 2025-12-14.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Prix Spot
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before the current column EL (the "01-oct." column).
$wsPrix.Range("EL1").EntireColumn.Insert()

# New header for the inserted column.
$wsPrix.Range("EL1").Value = "15-dec"

# The new column has no price data yet, mark every data row with "-".
$wsPrix.Range("EL2:EL25").Value = "-"

# ---------------------------------------------------------------------
# Sheet 2: Gaz
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A170").NumberFormat = "@"
$wsGaz.Range("A170").Value = "2025-12-13"
$wsGaz.Range("A170").ClearFormats()
$wsGaz.Range("B170").Value = 26.075

$wsGaz.Range("A171").NumberFormat = "@"
$wsGaz.Range("A171").Value = "2025-12-14"
$wsGaz.Range("A171").ClearFormats()
$wsGaz.Range("B171").Value = 26.075

# ---------------------------------------------------------------------
# Sheet 3: CO2
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A171").NumberFormat = "@"
$wsCo2.Range("A171").Value = "2025-12-13"
$wsCo2.Range("A171").ClearFormats()
$wsCo2.Range("B171").Value = 84.09999999999999

$wsCo2.Range("A172").NumberFormat = "@"
$wsCo2.Range("A172").Value = "2025-12-14"
$wsCo2.Range("A172").ClearFormats()
$wsCo2.Range("B172").Value = 84.09999999999999
